# "updae erik 14 nov" - corrections to payment rows on Sheet1 (Table1).
# Only true inputs are written here; formula-driven cells (I, K, L, H, N,
# and the Q192/N182 cached rollups) recompute automatically on recalc.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 211 (Eko / -) : payment date pushed a day, Pembayaran corrected ---
$ws.Range("A211").Value = 45528
$ws.Range("J211").Value = 200000

# --- Row 257 (Rohim) : period + tagihan corrected ---
$ws.Range("D257").Value = 45554
$ws.Range("E257").Value = 45573
$ws.Range("F257").Value = 460000

# --- Row 258 (Indra) : period start moved, tagihan now explicit & matches tertagih ---
$ws.Range("D258").Value = 45525
$ws.Range("D258").Interior.Color = 65535
$ws.Range("F258").Value = 3938000
$ws.Range("G258").Value = 3938000

# --- Row 260 (Purwadi) : tagihan now explicit & matches tertagih ---
$ws.Range("F260").Value = 1563000
$ws.Range("G260").Value = 1563000

# --- Row 261 (Perorangan-2 / Saiman) : tagihan now explicit & matches tertagih ---
$ws.Range("F261").Value = 361000
$ws.Range("G261").Value = 361000

# --- Row 262 (Perorangan-2 / Rizal) : tagihan now explicit & matches tertagih ---
$ws.Range("F262").Value = 515000
$ws.Range("G262").Value = 515000

# --- Row 264 (Nurdin) : Periode Mulai re-derived via the usual lookup formula,
#     tagihan now explicit & matches tertagih ---
$ws.Range("D264").Formula = '=IFERROR(LOOKUP(2,1/(($B$2:B263=B264)*($C$2:C263=C264)), $E$2:E263)+1, 0)'
$ws.Range("F264").Value = 197000

# --- Row 267 (Sunar) : tagihan now explicit & matches tertagih ---
$ws.Range("F267").Value = 1077000
$ws.Range("G267").Value = 1077000

# --- Row 269 (Purwadi) : tagihan now explicit & matches tertagih ---
$ws.Range("F269").Value = 1397000
$ws.Range("G269").Value = 1397000

# --- Row 270 (Indra) : tagihan now explicit & matches tertagih; Keterangan
#     overwritten with a literal note instead of the usual formula ---
$ws.Range("F270").Value = 6423000
$ws.Range("G270").Value = 6423000
$ws.Range("N270").Value = "kurang jay & robi"

# --- View state: scrolled down and a different cell selected ---
$excel.ActiveWindow.ScrollRow = 241
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K266").Select()
